# Updated symbol list on Sun Dec 25 03:46:42 UTC 2022 with GitHub Actions
#
# The "Price" column (D) holds numeric-looking values that are stored as
# plain text in the workbook (e.g. "245.65"), so every D-column write below
# first forces the cell to Text format ("@") to stop Excel's automatic
# number coercion from turning the literal string into a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
}

# --- simple price tweaks -------------------------------------------------
Set-TextValue "D2" "245.64"
Set-TextValue "D3" "23.04"
Set-TextValue "D4" "5.399"
Set-TextValue "D5" "0.06054"
Set-TextValue "D6" "3.393"
Set-TextValue "D7" "0.8067"
Set-TextValue "D8" "0.9312"

# --- rows 9-17: coin ranking reshuffled (each row shifts to the next slot,
#     values refreshed); row 17 receives what used to be row 9's coin ------
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D9" "0.1428"
$ws.Range("E9").Value = "8WazirXWRX"

$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D10" "0.07471"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"

$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D11" "0.03354"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.03070"
$ws.Range("E12").Value = "11BitrueCoinBTR"

$ws.Range("B13").Value = "MCDex"
$ws.Range("C13").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D13" "4.010"
$ws.Range("E13").Value = "12MCDexMCB"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09376"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D15" "0.001588"
$ws.Range("E15").Value = "14BitForexTokenBF"

$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D16" "0.04813"
$ws.Range("E16").Value = "15CoinExTokenCET"

$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D17" "0.0005942"
$ws.Range("E17").Value = "16OneONE"

# --- more simple price tweaks ---------------------------------------------
Set-TextValue "D18" "0.005252"
Set-TextValue "D20" "0.0009843"
Set-TextValue "D21" "0.00008704"
Set-TextValue "D23" "6.443"
Set-TextValue "D40" "0.03984"

# --- rows 41-43: another reshuffle ----------------------------------------
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1077"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002711"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003038"
$ws.Range("E43").Value = "42KickTokenKICK"

# --- trailing tweaks -------------------------------------------------------
Set-TextValue "D47" "0.0005802"
$ws.Range("E48").Value = "47CoinbaseStockTokenCOINBestin24h"
Set-TextValue "D49" "0.002198"
